$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-17 12:35:35"

# Drop the existing hyperlinks up front -- rows are about to be rewritten /
# reshuffled and the relationships would otherwise go stale.
$ws.Hyperlinks.Delete()

# --- Row 2: NEW - medical Teams/Graph API job (1st posting, detail/5415330) ---
$ws.Range("A2").Value = $newTimestamp
$ws.Range("B2").Value = "【医療機関向け】Teams連携「手術室予約承認システム」開発(Graph API/Azure)"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5415330"
$ws.Range("G2").Value = 265
$ws.Range("H2").Value = "🔥API ◆開発"

# --- Row 3: NEW - medical Teams/Graph API job (2nd posting, detail/5415235) ---
$ws.Range("A3").Value = $newTimestamp
$ws.Range("B3").Value = "【医療機関向け】Teams連携「手術室予約承認システム」開発(Graph API/Azure)"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5415235"
$ws.Range("G3").Value = 265
$ws.Range("H3").Value = "🔥API ◆開発"

# --- Row 4 (previously row 2): RPA tool job -- unchanged content, refreshed timestamp ---
$ws.Range("A4").Value = $newTimestamp
$ws.Range("B4").Value = "【募集】RPAツール「RoboTANGO」設定代行の専門家を探しています"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5405023"
$ws.Range("G4").Value = 178
$ws.Range("H4").Value = "★bot ◆ツール"

# --- Row 5 (previously row 3): Zoom consultation job -- unchanged content, refreshed timestamp ---
$ws.Range("A5").Value = $newTimestamp
$ws.Range("B5").Value = "【相談希望】在庫管理・出品補助ツールの開発に関するZoom面談依頼"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5398112"
$ws.Range("G5").Value = 158
$ws.Range("H5").Value = "◆ツール,開発 ◇管理"

# --- Row 6: NEW - KIntone customization job (detail/5415325) ---
$ws.Range("A6").Value = $newTimestamp
$ws.Range("B6").Value = "【急募】KIntoneアプリ間連携のカスタマイズ依頼"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5415325"
$ws.Range("G6").Value = 38
$ws.Range("H6").Value = "◇アプリ"

# --- Row 7 (previously row 4): pharmaceutical matching job -- unchanged content, refreshed timestamp ---
$ws.Range("A7").Value = $newTimestamp
$ws.Range("B7").Value = "【医薬品マッチング】高額医薬品の譲渡支援システム構築"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5415061"
$ws.Range("G7").Value = 33

# --- Row 8 (previously row 5): VBA quiz job -- unchanged content, refreshed timestamp ---
$ws.Range("A8").Value = $newTimestamp
$ws.Range("B8").Value = "初回 【急募・即決します】VBAで1問1答問題集の作成"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "~ 5,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5414812"
$ws.Range("G8").Value = 10

# --- Re-create the hyperlinks (and their visual style) for F2:F8, in row order ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5415330")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5415235")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5405023")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5398112")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5415325")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5415061")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5414812")
$ws.Range("F2:F8").Style = "Hyperlink"

# --- Column width adjustments (B: 38 -> 49, D: 28 -> 32) ---
$ws.Columns.Item(2).ColumnWidth = 49
$ws.Columns.Item(4).ColumnWidth = 32
